$wb = $excel.ActiveWorkbook

# Rename "ER" sheet to lowercase "er"
$erSheet = $wb.Worksheets.Item("ER")
$erSheet.Name = "er"

# Clear the explicit "applied number format" style from the data cells so
# they fall back to the default style (the diff removes the s="1" style
# index entirely, leaving only the default cellXfs entry).
$cellsToClear = @("C2","D2","E2","F2","J2","K2", `
                  "C3","D3","I3","J3","K3", `
                  "D7","E7","F7","I7","J7","K7", `
                  "C8","D8","F8","I8","K8","L8")

foreach ($addr in $cellsToClear) {
    $erSheet.Range($addr).ClearFormats()
}

# Update the active selection on the "er" sheet to H5
$erSheet.Range("H5").Select()

# Restore "TMA map" as the active sheet/tab (selection change above is scoped
# to the "er" sheet's own sheetView, but switching to it to set the selection
# would otherwise leave it as the active tab).
$wb.Worksheets.Item(1).Activate()
